$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The statement currently lists 5 workers across 14 rows (16-29). The
# update keeps only "JHOAN MAURICIO ALEMAN CUETO" and
# "MARIA JOSE MARRUGO PEÑATE", dropping the other three workers
# (JUAN SEBASTIAN RUIZ CUADROS, YESICA PAOLA CABARCAS SUAREZ and
# JULIO CESAR BENITEZ CORTEZ).

# Row 29 currently carries the special "closing" bottom-border formatting
# that marks the last row of the table. Row 19 (MARIA JOSE MARRUGO PEÑATE)
# will become the new last row once the rows below it are removed, so copy
# that closing format onto it first.
$ws.Range("B29:J29").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Drop the YESICA PAOLA CABARCAS SUAREZ and JULIO CESAR BENITEZ CORTEZ rows.
$ws.Range("20:29").EntireRow.Delete()

# Drop the JUAN SEBASTIAN RUIZ CUADROS rows.
$ws.Range("16:17").EntireRow.Delete()

# Refresh the summary figures at the top of the statement for the new,
# smaller set of workers/periods.
$ws.Range("E11").Value = 73333
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 1
